# "Screenshot is saved in same folder as the report and working fine"
#
# The underlying test-data change: the "click Opportunity record with
# value=..." step used by test cases TC103 (row 4) and TC104 (row 5) now
# points at the freshly-edited opportunity ("Tryedit") instead of the
# original one created earlier in the suite ("Try1"). TC105 (row 6) keeps
# referencing the original record and is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two cells whose step text now targets the edited record.
$ws.Range("J4").Value = 'click Opportunity record with value="Tryedit"'
$ws.Range("J5").Value = 'click Opportunity record with value="Tryedit"'

# Row heights were re-flowed (wrap-text rows resized) when the sheet was
# last saved.
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 105
$ws.Rows.Item(3).RowHeight = 135
$ws.Rows.Item(4).RowHeight = 120
$ws.Rows.Item(5).RowHeight = 105
$ws.Rows.Item(6).RowHeight = 105

# The sheet's scroll position/selection moved: the view now starts at row 3
# and the active cell is J6 instead of B6.
$ws.Range("J6").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
